$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new model row (row 7) for the added CSV run.
$ws.Range("A7").Value = "2_442"

# Convert the existing MODEL ID values (A2:A6) from numbers to the new
# text-based model identifiers used after the 1745 CSV file was added.
$ws.Range("A2").Value = "1_471"
$ws.Range("A3").Value = "2_471"
$ws.Range("A4").Value = "3_471"
$ws.Range("A5").Value = "4_471"
$ws.Range("A6").Value = "5_471"

$ws.Range("B7").Value = "Scaled Speed`nWeekdays oh`nDaypart oh"
$ws.Range("C7").Value = "LSTM(50)+dro(0.5)`nLSTM(50)+dro(0.5)`nLSTM(33)"
$ws.Range("D7").Value = 30
$ws.Range("E7").Value = "15m Back`n15m Forward"
$ws.Range("F7").Value = "1 WEEK"
$ws.Range("G7").Value = "Starting from`n16 of March, 8, up to end of May"
$ws.Range("H7").Value = "First 7 days of June"
$ws.Range("I7").Value = 15.5811
$ws.Range("J7").Value = 27.226
$ws.Range("K7").Value = 24.280687

# B/C/E/G/H carry the wrapped-text style used by the rest of the table.
$ws.Range("B7").WrapText = $true
$ws.Range("C7").WrapText = $true
$ws.Range("E7").WrapText = $true
$ws.Range("G7").WrapText = $true
$ws.Range("H7").WrapText = $true

$ws.Rows.Item(7).RowHeight = 56.25

$ws.Range("L7").Select()
